$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 29.285715
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = 93
$ws.Range("N2").Value = -376
$ws.Range("H12").Value = 3000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H40").Value = 2277.4
$ws.Range("I40").Value = 3285.7144
$ws.Range("J40").Value = 2063.5151
$ws.Range("K40").Value = 3285.7144
$ws.Range("L40").Value = 2063.5151
$ws.Range("M40").Value = -3110.7144
$ws.Range("N40").Value = -2413.5151
$ws.Range("H53").Value = 4039.1875
$ws.Range("I53").Value = 3962.875
$ws.Range("K53").Value = 3962.875
$ws.Range("M53").Value = -3325.875
$ws.Range("H70").Value = 83336840
$ws.Range("I70").Value = 41669096
$ws.Range("K70").Value = 125007288
$ws.Range("M70").Value = -125007018
$ws.Range("H73").Value = 83336840
$ws.Range("I73").Value = 41669096
$ws.Range("K73").Value = 125007288
$ws.Range("M73").Value = -125006352
$ws.Range("H99").Value = 386.5
$ws.Range("J99").Value = 200
$ws.Range("L99").Value = 600
$ws.Range("N99").Value = -3596
$ws.Range("H101").Value = 572.3333
$ws.Range("I101").Value = 309.57144
$ws.Range("K101").Value = 928.71432
$ws.Range("M101").Value = 693.28568
$ws.Range("H132").Value = 2365.7036
$ws.Range("I132").Value = 1911.7916
$ws.Range("K132").Value = 5735.3748
$ws.Range("M132").Value = -3205.3748
$ws.Range("H137").Value = 3261.7727
$ws.Range("I137").Value = 6050
$ws.Range("J137").Value = 2642.1667
$ws.Range("K137").Value = 18150
$ws.Range("L137").Value = 7926.500100000001
$ws.Range("M137").Value = -15600
$ws.Range("N137").Value = -13026.5001
$ws.Range("H141").Value = 2462.5454
$ws.Range("I141").Value = 1941.8572
$ws.Range("J141").Value = 3373.75
$ws.Range("K141").Value = 5825.571599999999
$ws.Range("L141").Value = 10121.25
$ws.Range("M141").Value = -645.5715999999993
$ws.Range("N141").Value = -20481.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1536.26
$ws.Range("I32").Value = 1494.433
$ws.Range("J32").Value = 2888.6667
$ws.Range("K32").Value = 1494.433
$ws.Range("L32").Value = 2888.6667
$ws.Range("M32").Value = -1207.433
$ws.Range("N32").Value = -3462.6667
$ws.Range("H35").Value = 2963
$ws.Range("I35").Value = 2963
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2963
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2557
$ws.Range("N35").ClearContents()
$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 5000
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = -4727
$ws.Range("H45").Value = 11483.857
$ws.Range("I45").Value = 1833.3334
$ws.Range("J45").Value = 18721.75
$ws.Range("K45").Value = 1833.3334
$ws.Range("L45").Value = 18721.75
$ws.Range("M45").Value = -1456.3334
$ws.Range("N45").Value = -19475.75
$ws.Range("H61").Value = 10037.956
$ws.Range("I61").Value = 3952
$ws.Range("K61").Value = 3952
$ws.Range("M61").Value = -3740
$ws.Range("H74").Value = 70529.21000000001
$ws.Range("I74").Value = 117234.36
$ws.Range("K74").Value = 117234.36
$ws.Range("M74").Value = -116360.36
$ws.Range("H77").Value = 70529.21000000001
$ws.Range("I77").Value = 117234.36
$ws.Range("K77").Value = 586171.8
$ws.Range("M77").Value = -581803.8
$ws.Range("H97").Value = 5567135.5
$ws.Range("I97").Value = 1163.8
$ws.Range("K97").Value = 1163.8
$ws.Range("M97").Value = -667.8
$ws.Range("H132").Value = 6373.9
$ws.Range("J132").Value = 9067.643
$ws.Range("L132").Value = 27202.929
$ws.Range("N132").Value = -32262.929
$ws.Range("H136").Value = 10037.956
$ws.Range("I136").Value = 3952
$ws.Range("K136").Value = 11856
$ws.Range("M136").Value = -9306

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4887
$ws.Range("H22").Value = 449.75
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -846
$ws.Range("H37").Value = 3149
$ws.Range("I37").Value = 1298
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 1298
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -1161
$ws.Range("N37").Value = -5274
$ws.Range("H86").Value = 32898138
$ws.Range("I86").Value = 13891626
$ws.Range("J86").Value = 50004000
$ws.Range("K86").Value = 13891626
$ws.Range("L86").Value = 50004000
$ws.Range("M86").Value = -13890503
$ws.Range("N86").Value = -50006246
$ws.Range("H89").Value = 32898138
$ws.Range("I89").Value = 13891626
$ws.Range("J89").Value = 50004000
$ws.Range("K89").Value = 69458130
$ws.Range("L89").Value = 250020000
$ws.Range("M89").Value = -69452514
$ws.Range("N89").Value = -250031232
$ws.Range("H99").Value = 3249427
$ws.Range("I99").Value = 2302.4443
$ws.Range("J99").Value = 9094251
$ws.Range("K99").Value = 2302.4443
$ws.Range("L99").Value = 9094251
$ws.Range("M99").Value = -804.4443000000001
$ws.Range("N99").Value = -9097247
$ws.Range("H134").Value = 6368.6313
$ws.Range("I134").Value = 2447.7273
$ws.Range("J134").Value = 11759.875
$ws.Range("K134").Value = 7343.1819
$ws.Range("L134").Value = 35279.625
$ws.Range("M134").Value = -4808.1819
$ws.Range("N134").Value = -40349.625

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9986.325000000001
$ws.Range("I31").Value = 4783.5557
$ws.Range("K31").Value = 4783.5557
$ws.Range("M31").Value = -4488.5557
$ws.Range("H34").Value = 9986.325000000001
$ws.Range("I34").Value = 4783.5557
$ws.Range("K34").Value = 4783.5557
$ws.Range("M34").Value = -4581.5557
$ws.Range("H105").Value = 3971226
$ws.Range("I105").Value = 5495570.5
$ws.Range("J105").Value = 7930
$ws.Range("K105").Value = 5495570.5
$ws.Range("L105").Value = 7930
$ws.Range("M105").Value = -5493823.5
$ws.Range("N105").Value = -11424
$ws.Range("H134").Value = 8079.115
$ws.Range("I134").Value = 1789.875
$ws.Range("J134").Value = 10874.333
$ws.Range("K134").Value = 5369.625
$ws.Range("L134").Value = 32622.999
$ws.Range("M134").Value = -2834.625
$ws.Range("N134").Value = -37692.999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 41670764
$ws.Range("I75").Value = 166671680
$ws.Range("J75").Value = 23813494
$ws.Range("K75").Value = 500015040
$ws.Range("L75").Value = 71440482
$ws.Range("M75").Value = -500014042
$ws.Range("N75").Value = -71442478
$ws.Range("H78").Value = 41670764
$ws.Range("I78").Value = 166671680
$ws.Range("J78").Value = 23813494
$ws.Range("K78").Value = 1500045120
$ws.Range("L78").Value = 214321446
$ws.Range("M78").Value = -1500040128
$ws.Range("N78").Value = -214331430
$ws.Range("H107").Value = 18182098
$ws.Range("J107").Value = 25000264
$ws.Range("L107").Value = 75000792
$ws.Range("N107").Value = -75004632
$ws.Range("H113").Value = 5126.4287
$ws.Range("I113").Value = 1343.8334
$ws.Range("J113").Value = 7963.375
$ws.Range("K113").Value = 4031.5002
$ws.Range("L113").Value = 23890.125
$ws.Range("M113").Value = -1861.5002
$ws.Range("N113").Value = -28230.125
$ws.Range("H122").Value = 2023331.1
$ws.Range("I122").Value = 4042127.5
$ws.Range("K122").Value = 36379147.5
$ws.Range("M122").Value = -36376697.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 200
$ws.Range("J10").Value = 200
$ws.Range("L10").Value = 200
$ws.Range("N10").Value = -538
$ws.Range("H80").Value = 3641.8
$ws.Range("I80").Value = 3064.182
$ws.Range("K80").Value = 3064.182
$ws.Range("M80").Value = -2066.182
$ws.Range("H83").Value = 3641.8
$ws.Range("I83").Value = 3064.182
$ws.Range("K83").Value = 15320.91
$ws.Range("M83").Value = -10328.91

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7999.25
$ws.Range("J16").Value = 7999.3335
$ws.Range("L16").Value = 7999.3335
$ws.Range("N16").Value = -8339.333500000001
$ws.Range("H46").Value = 2015.8
$ws.Range("I46").Value = 478.6
$ws.Range("J46").Value = 2784.4
$ws.Range("K46").Value = 478.6
$ws.Range("L46").Value = 2784.4
$ws.Range("M46").Value = -290.6
$ws.Range("N46").Value = -3160.4
$ws.Range("H82").Value = 63381760
$ws.Range("I82").Value = 101409450
$ws.Range("J82").Value = 2280.5
$ws.Range("K82").Value = 101409450
$ws.Range("L82").Value = 2280.5
$ws.Range("M82").Value = -101409089
$ws.Range("N82").Value = -3002.5
$ws.Range("H85").Value = 63381760
$ws.Range("I85").Value = 101409450
$ws.Range("J85").Value = 2280.5
$ws.Range("K85").Value = 101409450
$ws.Range("L85").Value = 2280.5
$ws.Range("M85").Value = -101408202
$ws.Range("N85").Value = -4776.5
$ws.Range("H93").Value = 4211.1787
$ws.Range("I93").Value = 3148.4285
$ws.Range("K93").Value = 3148.4285
$ws.Range("M93").Value = -1900.4285
$ws.Range("H100").Value = 5195.857
$ws.Range("I100").Value = 3289
$ws.Range("J100").Value = 6626
$ws.Range("K100").Value = 3289
$ws.Range("L100").Value = 6626
$ws.Range("M100").Value = -2748
$ws.Range("N100").Value = -7708

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3133
$ws.Range("I6").Value = 4000
$ws.Range("K6").Value = 4000
$ws.Range("M6").Value = -3885
$ws.Range("H125").Value = 50536
$ws.Range("J125").Value = 50536
$ws.Range("L125").Value = 50536
$ws.Range("N125").Value = -60376
$ws.Range("H132").Value = 14297556
$ws.Range("I132").Value = 23815260
$ws.Range("J132").Value = 20999.857
$ws.Range("K132").Value = 71445780
$ws.Range("L132").Value = 62999.571
$ws.Range("M132").Value = -71443250
$ws.Range("N132").Value = -68059.571
